# Auto-generated edit script: updates market-price-derived columns (H-N)
# for the Leve profit sheets, reflecting refreshed market data from the
# scheduled runner. Values for H, I, J come from the refreshed market
# snapshot; K = I*F, L = J*F, M = E-K (when K<>0), N = -(L+2E) (when L<>0).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1287.8235
$ws.Range("I70").Value = 950
$ws.Range("J70").Value = 1391.7693
$ws.Range("K70").Value = 2850
$ws.Range("L70").Value = 4175.3079
$ws.Range("M70").Value = -2580
$ws.Range("N70").Value = -4715.3079
$ws.Range("H73").Value = 1287.8235
$ws.Range("I73").Value = 950
$ws.Range("J73").Value = 1391.7693
$ws.Range("K73").Value = 2850
$ws.Range("L73").Value = 4175.3079
$ws.Range("M73").Value = -1914
$ws.Range("N73").Value = -6047.3079
$ws.Range("H74").Value = 2849.6667
$ws.Range("I74").Value = 2599.4285
$ws.Range("J74").Value = 3200
$ws.Range("K74").Value = 2599.4285
$ws.Range("L74").Value = 3200
$ws.Range("M74").Value = -1663.4285
$ws.Range("N74").Value = -5072
$ws.Range("H77").Value = 2849.6667
$ws.Range("I77").Value = 2599.4285
$ws.Range("J77").Value = 3200
$ws.Range("K77").Value = 12997.1425
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = -8317.1425
$ws.Range("N77").Value = -25360
$ws.Range("H80").Value = 3076.0833
$ws.Range("I80").Value = 2050.6667
$ws.Range("J80").Value = 4101.5
$ws.Range("K80").Value = 6152.000100000001
$ws.Range("L80").Value = 12304.5
$ws.Range("M80").Value = -5154.000100000001
$ws.Range("N80").Value = -14300.5
$ws.Range("H83").Value = 3076.0833
$ws.Range("I83").Value = 2050.6667
$ws.Range("J83").Value = 4101.5
$ws.Range("K83").Value = 18456.0003
$ws.Range("L83").Value = 36913.5
$ws.Range("M83").Value = -13464.0003
$ws.Range("N83").Value = -46897.5
$ws.Range("H132").Value = 24695348
$ws.Range("I132").Value = 3176510
$ws.Range("J132").Value = 100011290
$ws.Range("K132").Value = 9529530
$ws.Range("L132").Value = 300033870
$ws.Range("M132").Value = -9527000
$ws.Range("N132").Value = -300038930
$ws.Range("H135").Value = 30304268
$ws.Range("I135").Value = 1150.6
$ws.Range("J135").Value = 125001510
$ws.Range("K135").Value = 10355.4
$ws.Range("L135").Value = 1125013590
$ws.Range("M135").Value = -7820.4
$ws.Range("N135").Value = -1125018660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 18666.334
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 18666.334
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 18666.334
$ws.Range("N52").Value = -19302.334
$ws.Range("H61").Value = 1870.7037
$ws.Range("I61").Value = 1923.7693
$ws.Range("J61").Value = 1821.4286
$ws.Range("K61").Value = 1923.7693
$ws.Range("L61").Value = 1821.4286
$ws.Range("M61").Value = -1711.7693
$ws.Range("N61").Value = -2245.4286
$ws.Range("H136").Value = 1870.7037
$ws.Range("I136").Value = 1923.7693
$ws.Range("J136").Value = 1821.4286
$ws.Range("K136").Value = 5771.3079
$ws.Range("L136").Value = 5464.2858
$ws.Range("M136").Value = -3221.3079
$ws.Range("N136").Value = -10564.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1699.7407
$ws.Range("I86").Value = 1560.7778
$ws.Range("J86").Value = 1977.6666
$ws.Range("K86").Value = 1560.7778
$ws.Range("L86").Value = 1977.6666
$ws.Range("M86").Value = -437.7778000000001
$ws.Range("N86").Value = -4223.6666
$ws.Range("H89").Value = 1699.7407
$ws.Range("I89").Value = 1560.7778
$ws.Range("J89").Value = 1977.6666
$ws.Range("K89").Value = 7803.889
$ws.Range("L89").Value = 9888.333
$ws.Range("M89").Value = -2187.889
$ws.Range("N89").Value = -21120.333
$ws.Range("H105").Value = 2346.2917
$ws.Range("I105").Value = 2406.25
$ws.Range("J105").Value = 2316.3125
$ws.Range("K105").Value = 2406.25
$ws.Range("L105").Value = 2316.3125
$ws.Range("M105").Value = -659.25
$ws.Range("N105").Value = -5810.3125
$ws.Range("H134").Value = 558431.2
$ws.Range("I134").Value = 978708.3
$ws.Range("J134").Value = 2580.7742
$ws.Range("K134").Value = 2936124.9
$ws.Range("L134").Value = 7742.3226
$ws.Range("M134").Value = -2933589.9
$ws.Range("N134").Value = -12812.3226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2347.92
$ws.Range("I62").Value = 2310.4211
$ws.Range("J62").Value = 2466.6667
$ws.Range("K62").Value = 2310.4211
$ws.Range("L62").Value = 2466.6667
$ws.Range("M62").Value = -1686.4211
$ws.Range("N62").Value = -3714.6667
$ws.Range("H65").Value = 2347.92
$ws.Range("I65").Value = 2310.4211
$ws.Range("J65").Value = 2466.6667
$ws.Range("K65").Value = 11552.1055
$ws.Range("L65").Value = 12333.3335
$ws.Range("M65").Value = -8432.1055
$ws.Range("N65").Value = -18573.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1863.3684
$ws.Range("I109").Value = 920.7
$ws.Range("J109").Value = 2910.7778
$ws.Range("K109").Value = 2762.1
$ws.Range("L109").Value = 8732.3334
$ws.Range("M109").Value = -1722.1
$ws.Range("N109").Value = -10812.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29475
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29475
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29475
$ws.Range("N51").Value = -30493
$ws.Range("H57").Value = 14570.353
$ws.Range("I57").Value = 1500
$ws.Range("J57").Value = 15387.25
$ws.Range("K57").Value = 1500
$ws.Range("L57").Value = 15387.25
$ws.Range("M57").Value = -680
$ws.Range("N57").Value = -17027.25
$ws.Range("H80").Value = 3266.6667
$ws.Range("I80").Value = 2400
$ws.Range("J80").Value = 3700
$ws.Range("K80").Value = 2400
$ws.Range("L80").Value = 3700
$ws.Range("M80").Value = -1402
$ws.Range("N80").Value = -5696
$ws.Range("H83").Value = 3266.6667
$ws.Range("I83").Value = 2400
$ws.Range("J83").Value = 3700
$ws.Range("K83").Value = 12000
$ws.Range("L83").Value = 18500
$ws.Range("M83").Value = -7008
$ws.Range("N83").Value = -28484
$ws.Range("H123").Value = 12964.272
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 12964.272
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 12964.272
$ws.Range("N123").Value = -17864.272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 14388.889
$ws.Range("I68").Value = 51500
$ws.Range("J68").Value = 3785.7144
$ws.Range("K68").Value = 51500
$ws.Range("L68").Value = 3785.7144
$ws.Range("M68").Value = -50751
$ws.Range("N68").Value = -5283.7144
$ws.Range("H71").Value = 14388.889
$ws.Range("I71").Value = 51500
$ws.Range("J71").Value = 3785.7144
$ws.Range("K71").Value = 257500
$ws.Range("L71").Value = 18928.572
$ws.Range("M71").Value = -253756
$ws.Range("N71").Value = -26416.572
$ws.Range("H82").Value = 1360.1904
$ws.Range("I82").Value = 1666.2858
$ws.Range("J82").Value = 1207.1428
$ws.Range("K82").Value = 1666.2858
$ws.Range("L82").Value = 1207.1428
$ws.Range("M82").Value = -1305.2858
$ws.Range("N82").Value = -1929.1428
$ws.Range("H85").Value = 1360.1904
$ws.Range("I85").Value = 1666.2858
$ws.Range("J85").Value = 1207.1428
$ws.Range("K85").Value = 1666.2858
$ws.Range("L85").Value = 1207.1428
$ws.Range("M85").Value = -418.2858000000001
$ws.Range("N85").Value = -3703.1428
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1153.6
$ws.Range("I81").Value = 1167
$ws.Range("J81").Value = 1100
$ws.Range("K81").Value = 2334
$ws.Range("L81").Value = 2200
$ws.Range("M81").Value = -1273
$ws.Range("N81").Value = -4322
$ws.Range("H84").Value = 1153.6
$ws.Range("I84").Value = 1167
$ws.Range("J84").Value = 1100
$ws.Range("K84").Value = 11670
$ws.Range("L84").Value = 11000
$ws.Range("M84").Value = -6366
$ws.Range("N84").Value = -21608

